$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new rows 79-81 describing the new "corridas multiples" dashboard config section
$ws.Range("A79").Value = 10.1
$ws.Range("B79").Value = "dash_configs"
$ws.Range("C79").Value = "dash_configs"
$ws.Range("D79").Value = "alias_dash_lista"

$ws.Range("A80").Value = 10.2
$ws.Range("B80").Value = "dash_configs"
$ws.Range("C80").Value = "dash_configs"
$ws.Range("D80").Value = "alias_data_lista"

$ws.Range("A81").Value = 10.3
$ws.Range("B81").Value = "dash_configs"
$ws.Range("C81").Value = "dash_configs"
$ws.Range("D81").Value = "alias_insumos_lista"

# Apply same style (fill) as other "orden" column cells (style index 1 -> yellow fill)
$ws.Range("A79:A81").Interior.Color = 65535

# Re-sort the original table range by the "orden" column (matches the sortState
# left behind by Excel's Data > Sort, covering the pre-existing rows 2:78)
$dataRange = $ws.Range("A2:J78")
$sortKey = $ws.Range("A2:A78")
$dataRange.Sort($sortKey, 1, $null, $null, 1, $null, 1, 1)

# Update selection / view to match diff
$ws.Range("A6:XFD8").Select()
